$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 322-353 contain a duplicated copy of the "PAPA" crop block that
# should not be there (commit: "cambios en el cultivo de la papa (estab
# repetido)"). The real data that belongs in that slot is currently
# sitting in rows 354-385 ("PAPAYA" crop). Deleting rows 322-353 shifts
# the PAPAYA rows up into place and shrinks the used range to A1:J353.
$ws.Range("A322:J353").EntireRow.Delete()
